$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.448882
$ws.Cells.Item(2, 8).Value = 7.346646
$ws.Cells.Item(2, 9).Value = 0.3254773310529772
$ws.Cells.Item(2, 10).Value = 0.3254773310529772
$ws.Cells.Item(2, 13).Value = 1.123359666666667
$ws.Cells.Item(2, 14).Value = 3.370079
$ws.Cells.Item(2, 15).Value = 0.1748609861803909
$ws.Cells.Item(2, 16).Value = 0.1748609861803909
$ws.Cells.Item(2, 17).Value = 2.750975267226
$ws.Cells.Item(2, 18).Value = 24.758777405034
$ws.Cells.Item(2, 19).Value = 0.05691328708728516
$ws.Cells.Item(2, 20).Value = 0.05691328708728514

$ws.Cells.Item(3, 7).Value = 2.448882
$ws.Cells.Item(3, 8).Value = 7.346646
$ws.Cells.Item(3, 9).Value = 0.3254773310529772
$ws.Cells.Item(3, 10).Value = 0.3254773310529772
$ws.Cells.Item(3, 15).Value = 0.2317279222684864
$ws.Cells.Item(3, 16).Value = 0.2317279222684864
$ws.Cells.Item(3, 17).Value = 3.64562614458
$ws.Cells.Item(3, 18).Value = 32.81063530122
$ws.Cells.Item(3, 19).Value = 0.07542218567039871
$ws.Cells.Item(3, 20).Value = 0.0754221856703987

$ws.Cells.Item(4, 7).Value = 2.448882
$ws.Cells.Item(4, 8).Value = 7.346646
$ws.Cells.Item(4, 9).Value = 0.3254773310529772
$ws.Cells.Item(4, 10).Value = 0.3254773310529772
$ws.Cells.Item(4, 13).Value = 1.949056333333333
$ws.Cells.Item(4, 14).Value = 5.847169
$ws.Cells.Item(4, 15).Value = 0.3033880623283341
$ws.Cells.Item(4, 16).Value = 0.3033880623283341
$ws.Cells.Item(4, 17).Value = 4.773008971685999
$ws.Cells.Item(4, 18).Value = 42.957080745174
$ws.Cells.Item(4, 19).Value = 0.09874593679996048
$ws.Cells.Item(4, 20).Value = 0.09874593679996047

$ws.Cells.Item(5, 7).Value = 2.448882
$ws.Cells.Item(5, 8).Value = 7.346646
$ws.Cells.Item(5, 9).Value = 0.3254773310529772
$ws.Cells.Item(5, 10).Value = 0.3254773310529772
$ws.Cells.Item(5, 13).Value = 0.3313766666666667
$ws.Cells.Item(5, 14).Value = 0.9941300000000001
$ws.Cells.Item(5, 15).Value = 0.05158174398627213
$ws.Cells.Item(5, 16).Value = 0.05158174398627213
$ws.Cells.Item(5, 17).Value = 0.8115023542200001
$ws.Cells.Item(5, 18).Value = 7.30352118798
$ws.Cells.Item(5, 19).Value = 0.01678868836370981
$ws.Cells.Item(5, 20).Value = 0.0167886883637098

$ws.Cells.Item(6, 7).Value = 2.448882
$ws.Cells.Item(6, 8).Value = 7.346646
$ws.Cells.Item(6, 9).Value = 0.3254773310529772
$ws.Cells.Item(6, 10).Value = 0.3254773310529772
$ws.Cells.Item(6, 13).Value = 1.531818666666667
$ws.Cells.Item(6, 14).Value = 4.595456
$ws.Cells.Item(6, 15).Value = 0.2384412852365166
$ws.Cells.Item(6, 16).Value = 0.2384412852365165
$ws.Cells.Item(6, 17).Value = 3.751243160064
$ws.Cells.Item(6, 18).Value = 33.761188440576
$ws.Cells.Item(6, 19).Value = 0.07760723313162307
$ws.Cells.Item(6, 20).Value = 0.07760723313162304

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.075089
$ws.Cells.Item(7, 8).Value = 15.225267
$ws.Cells.Item(7, 9).Value = 0.6745226689470228
$ws.Cells.Item(7, 10).Value = 0.6745226689470227
$ws.Cells.Item(7, 13).Value = 1.123359666666667
$ws.Cells.Item(7, 14).Value = 3.370079
$ws.Cells.Item(7, 15).Value = 0.1748609861803909
$ws.Cells.Item(7, 16).Value = 0.1748609861803909
$ws.Cells.Item(7, 17).Value = 5.701150287343667
$ws.Cells.Item(7, 18).Value = 51.310352586093
$ws.Cells.Item(7, 19).Value = 0.1179476990931057
$ws.Cells.Item(7, 20).Value = 0.1179476990931057

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.075089
$ws.Cells.Item(8, 8).Value = 15.225267
$ws.Cells.Item(8, 9).Value = 0.6745226689470228
$ws.Cells.Item(8, 10).Value = 0.6745226689470227
$ws.Cells.Item(8, 15).Value = 0.2317279222684864
$ws.Cells.Item(8, 16).Value = 0.2317279222684864
$ws.Cells.Item(8, 17).Value = 7.555234243410001
$ws.Cells.Item(8, 18).Value = 67.99710819069
$ws.Cells.Item(8, 19).Value = 0.1563057365980877
$ws.Cells.Item(8, 20).Value = 0.1563057365980877

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.075089
$ws.Cells.Item(9, 8).Value = 15.225267
$ws.Cells.Item(9, 9).Value = 0.6745226689470228
$ws.Cells.Item(9, 10).Value = 0.6745226689470227
$ws.Cells.Item(9, 13).Value = 1.949056333333333
$ws.Cells.Item(9, 14).Value = 5.847169
$ws.Cells.Item(9, 15).Value = 0.3033880623283341
$ws.Cells.Item(9, 16).Value = 0.3033880623283341
$ws.Cells.Item(9, 17).Value = 9.891634357680333
$ws.Cells.Item(9, 18).Value = 89.02470921912301
$ws.Cells.Item(9, 19).Value = 0.2046421255283736
$ws.Cells.Item(9, 20).Value = 0.2046421255283736

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.075089
$ws.Cells.Item(10, 8).Value = 15.225267
$ws.Cells.Item(10, 9).Value = 0.6745226689470228
$ws.Cells.Item(10, 10).Value = 0.6745226689470227
$ws.Cells.Item(10, 13).Value = 0.3313766666666667
$ws.Cells.Item(10, 14).Value = 0.9941300000000001
$ws.Cells.Item(10, 15).Value = 0.05158174398627213
$ws.Cells.Item(10, 16).Value = 0.05158174398627213
$ws.Cells.Item(10, 17).Value = 1.681766075856667
$ws.Cells.Item(10, 18).Value = 15.13589468271
$ws.Cells.Item(10, 19).Value = 0.03479305562256232
$ws.Cells.Item(10, 20).Value = 0.03479305562256231

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 5.075089
$ws.Cells.Item(11, 8).Value = 15.225267
$ws.Cells.Item(11, 9).Value = 0.6745226689470228
$ws.Cells.Item(11, 10).Value = 0.6745226689470227
$ws.Cells.Item(11, 13).Value = 1.531818666666667
$ws.Cells.Item(11, 14).Value = 4.595456
$ws.Cells.Item(11, 15).Value = 0.2384412852365166
$ws.Cells.Item(11, 16).Value = 0.2384412852365165
$ws.Cells.Item(11, 17).Value = 7.774116065194668
$ws.Cells.Item(11, 18).Value = 69.96704458675201
$ws.Cells.Item(11, 19).Value = 0.1608340521048935
$ws.Cells.Item(11, 20).Value = 0.1608340521048935
